# Update to TrancheDetails Sheet
# - Rename the "TrancheBounds" sheet to "TrancheDetails"
# - Rename the "Mean Gas Rate" row label to "Normalised Gas Rate"
# - Select the full column A (as if the column header were clicked)
# - Widen column A to fit the new, longer label

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("TrancheBounds")
$ws.Name = "TrancheDetails"

$ws.Range("A4").Value = "Normalised Gas Rate"

$ws.Columns.Item(1).ColumnWidth = 17.46

$ws.Columns.Item(1).Select()
